# Update Reaction_number values (column C) on the NBR and BAR sheets.

$wb = $excel.ActiveWorkbook

$nbr = $wb.Worksheets.Item("NBR")
$bar = $wb.Worksheets.Item("BAR")

$nbrValues = @(843, 831, 824, 822, 811, 0, 799, 797, 798, 758, 0, 0, 743, 742, 727, 717, 711, 710, 701)
$barValues = @(676, 678, 684, 686, 686, 0, 688, 683, 681, 681, 0, 0, 679, 663, 661, 660, 660, 660, 659)

for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $nbr.Cells.Item($row, 3).Value = $nbrValues[$i]
}

for ($i = 0; $i -lt $barValues.Length; $i++) {
    $row = $i + 2
    $bar.Cells.Item($row, 3).Value = $barValues[$i]
}
